# "Generate Report for Archive"
# The localization status has moved on from "Ready for handoff" to
# "In Translation" everywhere that string is shown (the Overview rollup
# columns for zh-cn / de-de, and the per-language Status column), and the
# now-shorter status text means the two status columns can be narrower.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# New (narrower) width for the status-ish columns, expressed the way Excel's
# COM ColumnWidth property takes it (character units at the workbook's
# default font). 12.5 lands in the same rounded bucket the target width maps
# to.
$newColWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    $rowOffset = $used.Row
    $colOffset = $used.Column

    for ($r = 0; $r -lt $rows; $r++) {
        for ($c = 0; $c -lt $cols; $c++) {
            $cell = $ws.Cells.Item($rowOffset + $r, $colOffset + $c)
            # NB: put the string literal on the left of -eq so PowerShell
            # compares as a string instead of coercing against a boolean
            # cell value (e.g. the "True"/"False" cells elsewhere on these
            # sheets).
            if ($oldStatus -eq $cell.Value2) {
                $cell.Value = $newStatus
            }
        }
    }
}

# Narrow the two "status" columns on the Overview sheet (zh-cn / de-de,
# columns E and F) ...
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newColWidth
$overview.Columns.Item(6).ColumnWidth = $newColWidth

# ... and the Status column (C) on each per-language report sheet.
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newColWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newColWidth
